$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert 8 new blank rows above the existing header row, pushing the
#        old row 1 (Account Name / Date / ... / Bad Order) down to row 9.
$ws.Range("A1:A8").EntireRow.Insert()

# --- 2. Row 1: big merged title placeholder (bold, 20pt, centered), blank.
$r1 = $ws.Range("A1:G1")
$r1.Font.Bold = $true
$r1.Font.Size = 20
$r1.HorizontalAlignment = -4108
$r1.RowHeight = 26.25
$r1.Merge()

# --- 3. Rows 2-5: merged placeholder rows (italic, centered), blank.
foreach ($r in 2..5) {
    $rng = $ws.Range("A" + $r + ":G" + $r)
    $rng.Font.Italic = $true
    $rng.HorizontalAlignment = -4108
    $rng.Merge()
}

# --- 4. Row 7: "Collection Register" title, merged, bold (non-italic) + plain
#        centered companions.
$ws.Range("A7").Value = "Collection Register"
$a7 = $ws.Range("A7")
$a7.Font.Bold = $true
$a7.Font.Italic = $false
$a7.Font.Size = 11
$a7.HorizontalAlignment = -4108

$rest7 = $ws.Range("B7:G7")
$rest7.Font.Bold = $false
$rest7.Font.Italic = $false
$rest7.Font.Size = 11
$rest7.HorizontalAlignment = -4108

$ws.Range("A7:G7").Merge()

# --- 5. Row 8: thin spacer row just before the table header.
$ws.Rows.Item(8).RowHeight = 15.75

# --- 6. Row 9 (former row 1): table header - add a thick outline box (medium
#        top/bottom across, medium left on A9, medium right on G9) and keep
#        the existing bold/centered header formatting + row height.
$ws.Rows.Item(9).RowHeight = 15.75

$hdr = $ws.Range("A9:G9")
$hdr.Borders.Item(9).LineStyle = 1   # xlEdgeBottom
$hdr.Borders.Item(9).Weight = 4      # xlMedium
$hdr.Borders.Item(8).LineStyle = 1   # xlEdgeTop
$hdr.Borders.Item(8).Weight = 4      # xlMedium

$ws.Range("A9").Borders.Item(7).LineStyle = 1   # xlEdgeLeft
$ws.Range("A9").Borders.Item(7).Weight = 4

$ws.Range("G9").Borders.Item(10).LineStyle = 1  # xlEdgeRight
$ws.Range("G9").Borders.Item(10).Weight = 4

# --- 7. Move the active selection where Excel left it after these edits.
$ws.Range("C10").Select()

Write-Output "done"
